$app = $ppt
Write-Host ($app.Presentations | Get-Member | Out-String)
